$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Find_Customer" worksheet right after "Customer".
# ---------------------------------------------------------------------------
$customer = $wb.Worksheets.Item("Customer")
$find = $wb.Worksheets.Add($null, $customer)
$find.Name = "Find_Customer"

# ---------------------------------------------------------------------------
# 2. Populate the top of the new sheet first (this controls the order new
#    strings are interned in the shared-string table).
# ---------------------------------------------------------------------------
$find.Range("A1").Value = "What do you want to find?"
$find.Range("B1").Value = "Your input"

$find.Range("A2").Value = "Phone"
$find.Range("B2").Value = "'0328150801"

# ---------------------------------------------------------------------------
# 3. Update the Customer sheet's Phone column: add a leading zero to each
#    phone number (328150801 -> 0328150801, etc.), keeping them as text.
# ---------------------------------------------------------------------------
$customer.Range("E2").Value = "'0328150801"
$customer.Range("E3").Value = "'0328150802"
$customer.Range("E4").Value = "'0328150803"
$customer.Range("E5").Value = "'0328150804"
$customer.Range("E6").Value = "'0328150805"
$customer.Range("E7").Value = "'0328150806"

# Move the stored selection on the Customer sheet (it is no longer the
# active tab once Find_Customer is inserted/activated).
$customer.Range("G22").Select()

# ---------------------------------------------------------------------------
# 4. Finish populating the Find_Customer helper layout.
# ---------------------------------------------------------------------------
$find.Range("A3").Value = "Email"
$find.Range("B3").Value = "b@gmail.com"
$find.Hyperlinks.Add($find.Range("B3"), "mailto:b@gmail.com", "", "", "b@gmail.com")
$find.Range("B3").Style = "Hyperlink"

$find.Range("A4").Value = "FirstName"
$find.Range("B4").Value = "L"

$find.Range("A5").Value = "LastName"
$find.Range("B5").Value = "Nguyen"

# Column widths roughly matching the authored layout.
$find.Columns.Item(1).ColumnWidth = 28
$find.Columns.Item(2).ColumnWidth = 23.5
$find.Columns.Item(3).ColumnWidth = 16.333333333333332

# Active cell / selection on the new (now active) sheet.
$find.Range("B4").Select()
